# GradeSheet for 1 to 10 Fixed nearly
#
# Class 12 rows for "Compulsory Nepali" / "Social Studies & Life Skills"
# did not match the equivalent class 11 rows, so bring them in line:
#   - row 15 ("Compulsory Nepali", class 12) -> faculty should read
#     "Compulsory", same as row 3 ("Compulsory Nepali", class 11).
#   - row 16 ("Social Studies & Life Skills", class 12) -> is_compulsory
#     should be 0, same as row 4 ("Social Studies & Life Skills", class 11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K15").Value = "Compulsory"
$ws.Range("L16").Value = 0

# Minor formatting touch-ups that came along with the fix.
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Columns.Item(12).ColumnWidth = 13.140625

# Leave the selection on row 4, as it was when the fix was made.
$ws.Rows.Item(4).Select() | Out-Null
